# Auto-applied numeric corrections to H:N profit columns across multiple sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 125
$ws.Range("H125").Value = 8277.333000000001
$ws.Range("I125").Value = 23032
$ws.Range("J125").Value = 900
$ws.Range("K125").Value = 207288
$ws.Range("L125").Value = 8100
$ws.Range("M125").Value = -204828
$ws.Range("N125").Value = -13020

# row 129
$ws.Range("H129").Value = 1055.9762
$ws.Range("I129").Value = 738.1429000000001
$ws.Range("J129").Value = 1119.5428
$ws.Range("K129").Value = 2214.4287
$ws.Range("L129").Value = 3358.6284
$ws.Range("M129").Value = 2785.5713
$ws.Range("N129").Value = -13358.6284

# row 131
$ws.Range("H131").Value = 3227.5
$ws.Range("I131").Value = 1975.7142
$ws.Range("J131").Value = 4980
$ws.Range("K131").Value = 5927.142599999999
$ws.Range("L131").Value = 14940
$ws.Range("M131").Value = -887.1425999999992
$ws.Range("N131").Value = -25020

# row 135
$ws.Range("H135").Value = 1928.5676
$ws.Range("I135").Value = 1908.7354
$ws.Range("J135").Value = 2153.3333
$ws.Range("K135").Value = 17178.6186
$ws.Range("L135").Value = 19379.9997
$ws.Range("M135").Value = -14643.6186
$ws.Range("N135").Value = -24449.9997

# row 137
$ws.Range("H137").Value = 1473.7567
$ws.Range("I137").Value = 964.3684
$ws.Range("J137").Value = 2011.4445
$ws.Range("K137").Value = 2893.1052
$ws.Range("L137").Value = 6034.333500000001
$ws.Range("M137").Value = -343.1052
$ws.Range("N137").Value = -11134.3335

# row 138
$ws.Range("H138").Value = 4162.375
$ws.Range("I138").Value = 1272.8611
$ws.Range("J138").Value = 6162.8076
$ws.Range("K138").Value = 3818.5833
$ws.Range("L138").Value = 18488.4228
$ws.Range("M138").Value = 1321.4167
$ws.Range("N138").Value = -28768.4228

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3200.611
$ws.Range("I32").Value = 2736.0557
$ws.Range("K32").Value = 2736.0557
$ws.Range("M32").Value = -2449.0557

# row 61
$ws.Range("H61").Value = 2101.1667
$ws.Range("I61").Value = 1890.9375
$ws.Range("J61").Value = 2773.9
$ws.Range("K61").Value = 1890.9375
$ws.Range("L61").Value = 2773.9
$ws.Range("M61").Value = -1678.9375
$ws.Range("N61").Value = -3197.9

# row 74
$ws.Range("H74").Value = 1679.2703
$ws.Range("I74").Value = 1375
$ws.Range("J74").Value = 2179.1428
$ws.Range("K74").Value = 1375
$ws.Range("L74").Value = 2179.1428
$ws.Range("M74").Value = -501
$ws.Range("N74").Value = -3927.1428

# row 77
$ws.Range("H77").Value = 1679.2703
$ws.Range("I77").Value = 1375
$ws.Range("J77").Value = 2179.1428
$ws.Range("K77").Value = 6875
$ws.Range("L77").Value = 10895.714
$ws.Range("M77").Value = -2507
$ws.Range("N77").Value = -19631.714

# row 102
$ws.Range("H102").Value = 3377851.5
$ws.Range("I102").Value = 3715435.5
$ws.Range("J102").Value = 2011
$ws.Range("K102").Value = 3715435.5
$ws.Range("L102").Value = 2011
$ws.Range("M102").Value = -3713813.5
$ws.Range("N102").Value = -5255

# row 122
$ws.Range("H122").Value = 8548608
$ws.Range("I122").Value = 8548608
$ws.Range("K122").Value = 25645824
$ws.Range("M122").Value = -25643374

# row 132
$ws.Range("H132").Value = 2538.2188
$ws.Range("I132").Value = 2424.84
$ws.Range("K132").Value = 7274.52
$ws.Range("M132").Value = -4744.52

# row 136
$ws.Range("H136").Value = 2101.1667
$ws.Range("I136").Value = 1890.9375
$ws.Range("J136").Value = 2773.9
$ws.Range("K136").Value = 5672.8125
$ws.Range("L136").Value = 8321.700000000001
$ws.Range("M136").Value = -3122.8125
$ws.Range("N136").Value = -13421.7

$ws = $wb.Worksheets.Item("BSM")
# row 94
$ws.Range("H94").Value = 1948.36
$ws.Range("I94").Value = 1343.8462
$ws.Range("J94").Value = 2603.25
$ws.Range("K94").Value = 1343.8462
$ws.Range("L94").Value = 2603.25
$ws.Range("M94").Value = -892.8462
$ws.Range("N94").Value = -3505.25

# row 107
$ws.Range("H107").Value = 1636.4706
$ws.Range("I107").Value = 1378.4615
$ws.Range("J107").Value = 2475
$ws.Range("K107").Value = 1378.4615
$ws.Range("L107").Value = 2475
$ws.Range("M107").Value = 541.5385000000001
$ws.Range("N107").Value = -6315

# row 134
$ws.Range("H134").Value = 3024.842
$ws.Range("I134").Value = 2585.7778
$ws.Range("J134").Value = 3420
$ws.Range("K134").Value = 7757.3334
$ws.Range("L134").Value = 10260
$ws.Range("M134").Value = -5222.3334
$ws.Range("N134").Value = -15330

$ws = $wb.Worksheets.Item("CRP")
# row 11
$ws.Range("H11").Value = 1005
$ws.Range("I11").Value = 1005
$ws.Range("K11").Value = 1005
$ws.Range("M11").Value = -865

# row 16
$ws.Range("H16").Value = 3302.5
$ws.Range("I16").Value = 1833.3334
$ws.Range("J16").Value = 3932.1428
$ws.Range("K16").Value = 1833.3334
$ws.Range("L16").Value = 3932.1428
$ws.Range("M16").Value = -1546.3334
$ws.Range("N16").Value = -4506.1428

# row 58
$ws.Range("H58").Value = 1474.25
$ws.Range("I58").Value = 1316.7646
$ws.Range("J58").Value = 2366.6667
$ws.Range("K58").Value = 1316.7646
$ws.Range("L58").Value = 2366.6667
$ws.Range("M58").Value = -1113.7646
$ws.Range("N58").Value = -2772.6667

# row 63
$ws.Range("H63").Value = 42300
$ws.Range("J63").Value = 42300
$ws.Range("L63").Value = 42300
$ws.Range("N63").Value = -43672

# row 66
$ws.Range("H66").Value = 42300
$ws.Range("J66").Value = 42300
$ws.Range("L66").Value = 126900
$ws.Range("N66").Value = -133764

# row 105
$ws.Range("H105").Value = 2456.2856
$ws.Range("I105").Value = 2077.1428
$ws.Range("J105").Value = 2835.4285
$ws.Range("K105").Value = 2077.1428
$ws.Range("L105").Value = 2835.4285
$ws.Range("M105").Value = -330.1428000000001
$ws.Range("N105").Value = -6329.4285

# row 113
$ws.Range("H113").Value = 3302.5
$ws.Range("I113").Value = 1833.3334
$ws.Range("J113").Value = 3932.1428
$ws.Range("K113").Value = 1833.3334
$ws.Range("L113").Value = 3932.1428
$ws.Range("M113").Value = 336.6666
$ws.Range("N113").Value = -8272.1428

# row 136
$ws.Range("H136").Value = 1474.25
$ws.Range("I136").Value = 1316.7646
$ws.Range("J136").Value = 2366.6667
$ws.Range("K136").Value = 3950.2938
$ws.Range("L136").Value = 7100.000100000001
$ws.Range("M136").Value = -1400.2938
$ws.Range("N136").Value = -12200.0001

$ws = $wb.Worksheets.Item("CUL")
# row 23
$ws.Range("H23").Value = 12500086
$ws.Range("I23").Value = 20000020
$ws.Range("J23").Value = 196.66667
$ws.Range("K23").Value = 60000060
$ws.Range("L23").Value = 590.00001
$ws.Range("M23").Value = -59999825
$ws.Range("N23").Value = -1060.00001

# row 68
$ws.Range("H68").Value = 3227.125
$ws.Range("I68").Value = 4411.357
$ws.Range("J68").Value = 1569.2
$ws.Range("K68").Value = 13234.071
$ws.Range("L68").Value = 4707.6
$ws.Range("M68").Value = -12423.071
$ws.Range("N68").Value = -6329.6

# row 71
$ws.Range("H71").Value = 3227.125
$ws.Range("I71").Value = 4411.357
$ws.Range("J71").Value = 1569.2
$ws.Range("K71").Value = 39702.213
$ws.Range("L71").Value = 14122.8
$ws.Range("M71").Value = -35646.213
$ws.Range("N71").Value = -22234.8

# row 113
$ws.Range("H113").Value = 1395829.9
$ws.Range("I113").Value = 2000456.4
$ws.Range("J113").Value = 556070.75
$ws.Range("K113").Value = 6001369.199999999
$ws.Range("L113").Value = 1668212.25
$ws.Range("M113").Value = -5999199.199999999
$ws.Range("N113").Value = -1672552.25

# row 131
$ws.Range("H131").Value = 21154930
$ws.Range("J131").Value = 28572754
$ws.Range("L131").Value = 85718262
$ws.Range("N131").Value = -85728342

# row 136
$ws.Range("H136").Value = 12472.2
$ws.Range("I136").Value = 17287
$ws.Range("J136").Value = 5250
$ws.Range("K136").Value = 51861
$ws.Range("L136").Value = 15750
$ws.Range("M136").Value = -46761
$ws.Range("N136").Value = -25950

$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 5676.1606
$ws.Range("I70").Value = 5726.364
$ws.Range("J70").Value = 5492.0835
$ws.Range("K70").Value = 5726.364
$ws.Range("L70").Value = 5492.0835
$ws.Range("M70").Value = -5456.364
$ws.Range("N70").Value = -6032.0835

# row 73
$ws.Range("H73").Value = 5676.1606
$ws.Range("I73").Value = 5726.364
$ws.Range("J73").Value = 5492.0835
$ws.Range("K73").Value = 5726.364
$ws.Range("L73").Value = 5492.0835
$ws.Range("M73").Value = -4790.364
$ws.Range("N73").Value = -7364.0835

# row 122
$ws.Range("H122").Value = 3088702.2
$ws.Range("I122").Value = 4323183.5
$ws.Range("J122").Value = 2499.6667
$ws.Range("K122").Value = 12969550.5
$ws.Range("L122").Value = 7499.000100000001
$ws.Range("M122").Value = -12967100.5
$ws.Range("N122").Value = -12399.0001

$ws = $wb.Worksheets.Item("LTW")
# row 46
$ws.Range("H46").Value = 22223942
$ws.Range("I46").Value = 47620240
$ws.Range("J46").Value = 2178
$ws.Range("K46").Value = 47620240
$ws.Range("L46").Value = 2178
$ws.Range("M46").Value = -47620052
$ws.Range("N46").Value = -2554

# row 122
$ws.Range("H122").Value = 5091574
$ws.Range("I122").Value = 5496875.5
$ws.Range("J122").Value = 3335266.8
$ws.Range("K122").Value = 16490626.5
$ws.Range("L122").Value = 10005800.4
$ws.Range("M122").Value = -16488176.5
$ws.Range("N122").Value = -10010700.4

# row 132
$ws.Range("H132").Value = 16053967
$ws.Range("I132").Value = 18845254
$ws.Range("J132").Value = 4062.5
$ws.Range("K132").Value = 56535762
$ws.Range("L132").Value = 12187.5
$ws.Range("M132").Value = -56533232
$ws.Range("N132").Value = -17247.5

$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 15640.8
$ws.Range("I62").Value = 4250
$ws.Range("J62").Value = 23234.666
$ws.Range("K62").Value = 4250
$ws.Range("L62").Value = 23234.666
$ws.Range("M62").Value = -3626
$ws.Range("N62").Value = -24482.666

# row 65
$ws.Range("H65").Value = 15640.8
$ws.Range("I65").Value = 4250
$ws.Range("J65").Value = 23234.666
$ws.Range("K65").Value = 21250
$ws.Range("L65").Value = 116173.33
$ws.Range("M65").Value = -18130
$ws.Range("N65").Value = -122413.33

# row 96
$ws.Range("H96").Value = 1750
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = -627
$ws.Range("N96").Value = -3746

# row 107
$ws.Range("H107").Value = 58824684
$ws.Range("I107").Value = 90910460
$ws.Range("J107").Value = 753.3333
$ws.Range("K107").Value = 272731380
$ws.Range("L107").Value = 2259.9999
$ws.Range("M107").Value = -272729460
$ws.Range("N107").Value = -6099.9999

# row 122
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 12000
$ws.Range("N122").Value = -16900
$ws.Range("M122").ClearContents()

# row 126
$ws.Range("H126").Value = 1527.2727
$ws.Range("I126").Value = 1132.6666
$ws.Range("J126").Value = 2000.8
$ws.Range("K126").Value = 3397.9998
$ws.Range("L126").Value = 6002.4
$ws.Range("M126").Value = -927.9998000000001
$ws.Range("N126").Value = -10942.4

# row 136
$ws.Range("H136").Value = 1258.3334
$ws.Range("I136").Value = 1020
$ws.Range("J136").Value = 2450
$ws.Range("K136").Value = 3060
$ws.Range("L136").Value = 7350
$ws.Range("M136").Value = -510
$ws.Range("N136").Value = -12450
